# Update odds/volume figures in the "Jogos do Dia" sheet (rows 2-14) to match
# the latest Betfair Back/Lay snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 1.39
$ws.Cells.Item(2, 14).Value = 3.45
$ws.Cells.Item(2, 15).Value = 1.38
$ws.Cells.Item(2, 16).Value = 1.83
$ws.Cells.Item(2, 17).Value = 2.14
$ws.Cells.Item(2, 19).Value = 4
$ws.Cells.Item(2, 24).Value = 14
$ws.Cells.Item(2, 25).Value = 26
$ws.Cells.Item(2, 36).Value = 11
$ws.Cells.Item(2, 39).Value = 400
$ws.Cells.Item(3, 7).Value = 2.9
$ws.Cells.Item(3, 14).Value = 3.75
$ws.Cells.Item(3, 15).Value = 1.35
$ws.Cells.Item(3, 16).Value = 1.91
$ws.Cells.Item(3, 17).Value = 2.04
$ws.Cells.Item(3, 19).Value = 3.65
$ws.Cells.Item(3, 20).Value = 1.8
$ws.Cells.Item(3, 21).Value = 2.2
$ws.Cells.Item(3, 24).Value = 13.5
$ws.Cells.Item(3, 25).Value = 11
$ws.Cells.Item(3, 26).Value = 18.5
$ws.Cells.Item(3, 27).Value = 100
$ws.Cells.Item(3, 28).Value = 11.5
$ws.Cells.Item(3, 30).Value = 12.5
$ws.Cells.Item(3, 31).Value = 32
$ws.Cells.Item(3, 32).Value = 19
$ws.Cells.Item(3, 33).Value = 13
$ws.Cells.Item(3, 34).Value = 17
$ws.Cells.Item(3, 35).Value = 44
$ws.Cells.Item(3, 36).Value = 44
$ws.Cells.Item(3, 37).Value = 32
$ws.Cells.Item(3, 38).Value = 44
$ws.Cells.Item(3, 39).Value = 110
$ws.Cells.Item(3, 40).Value = 29
$ws.Cells.Item(3, 41).Value = 28
$ws.Cells.Item(4, 6).Value = 1.42
$ws.Cells.Item(4, 9).Value = 8.4
$ws.Cells.Item(4, 10).Value = 5.6
$ws.Cells.Item(4, 11).Value = 5.7
$ws.Cells.Item(4, 15).Value = 1.15
$ws.Cells.Item(4, 17).Value = 1.47
$ws.Cells.Item(4, 19).Value = 2.16
$ws.Cells.Item(4, 20).Value = 1.69
$ws.Cells.Item(4, 21).Value = 2.36
$ws.Cells.Item(4, 24).Value = 34
$ws.Cells.Item(4, 25).Value = 40
$ws.Cells.Item(4, 26).Value = 80
$ws.Cells.Item(4, 27).Value = 260
$ws.Cells.Item(4, 28).Value = 14
$ws.Cells.Item(4, 30).Value = 32
$ws.Cells.Item(4, 31).Value = 95
$ws.Cells.Item(4, 35).Value = 80
$ws.Cells.Item(4, 39).Value = 80
$ws.Cells.Item(4, 41).Value = 85
$ws.Cells.Item(5, 6).Value = 2.18
$ws.Cells.Item(5, 7).Value = 2.22
$ws.Cells.Item(5, 18).Value = 1.47
$ws.Cells.Item(5, 24).Value = 19
$ws.Cells.Item(5, 32).Value = 15.5
$ws.Cells.Item(5, 35).Value = 50
$ws.Cells.Item(5, 36).Value = 28
$ws.Cells.Item(5, 39).Value = 1000
$ws.Cells.Item(5, 40).Value = 13.5
$ws.Cells.Item(6, 16).Value = 1.8
$ws.Cells.Item(6, 21).Value = 1.92
$ws.Cells.Item(6, 39).Value = 140
$ws.Cells.Item(7, 6).Value = 3.25
$ws.Cells.Item(7, 7).Value = 3.35
$ws.Cells.Item(7, 8).Value = 2.42
$ws.Cells.Item(7, 20).Value = 1.78
$ws.Cells.Item(7, 21).Value = 2.18
$ws.Cells.Item(7, 31).Value = 26
$ws.Cells.Item(8, 6).Value = 1.89
$ws.Cells.Item(8, 8).Value = 4.7
$ws.Cells.Item(8, 9).Value = 4.9
$ws.Cells.Item(8, 21).Value = 1.99
$ws.Cells.Item(8, 34).Value = 22
$ws.Cells.Item(9, 16).Value = 2.1
$ws.Cells.Item(9, 25).Value = 11
$ws.Cells.Item(9, 26).Value = 15.5
$ws.Cells.Item(9, 31).Value = 21
$ws.Cells.Item(9, 35).Value = 32
$ws.Cells.Item(10, 9).Value = 3.55
$ws.Cells.Item(10, 15).Value = 1.45
$ws.Cells.Item(10, 21).Value = 1.94
$ws.Cells.Item(10, 25).Value = 11
$ws.Cells.Item(10, 30).Value = 15.5
$ws.Cells.Item(10, 32).Value = 14
$ws.Cells.Item(11, 7).Value = 2.4
$ws.Cells.Item(11, 10).Value = 3.2
$ws.Cells.Item(11, 11).Value = 3.25
$ws.Cells.Item(11, 24).Value = 10
$ws.Cells.Item(11, 32).Value = 14
$ws.Cells.Item(11, 40).Value = 28
$ws.Cells.Item(12, 6).Value = 8.800000000000001
$ws.Cells.Item(12, 7).Value = 9.4
$ws.Cells.Item(12, 10).Value = 5.1
$ws.Cells.Item(12, 11).Value = 5.3
$ws.Cells.Item(12, 18).Value = 1.56
$ws.Cells.Item(12, 20).Value = 1.93
$ws.Cells.Item(12, 21).Value = 2
$ws.Cells.Item(13, 7).Value = 1.74
$ws.Cells.Item(13, 17).Value = 1.89
$ws.Cells.Item(13, 29).Value = 9
$ws.Cells.Item(14, 8).Value = 1.78
$ws.Cells.Item(14, 9).Value = 1.79
$ws.Cells.Item(14, 11).Value = 4.2
$ws.Cells.Item(14, 14).Value = 4.4
$ws.Cells.Item(14, 16).Value = 2.16
$ws.Cells.Item(14, 17).Value = 1.84
$ws.Cells.Item(14, 20).Value = 1.8
$ws.Cells.Item(14, 21).Value = 2.18
$ws.Cells.Item(14, 24).Value = 19
$ws.Cells.Item(14, 25).Value = 9.800000000000001
$ws.Cells.Item(14, 26).Value = 10.5
$ws.Cells.Item(14, 39).Value = 90
